$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 22), shrinking the used range to A1:H21
$ws.Rows.Item(22).Delete()

# Shift the sensor columns (C:H) for rows 3-21 down from what's currently in rows 2-20,
# working from the bottom up so we don't clobber values before reading them.
for ($r = 21; $r -ge 3; $r--) {
    $src = $r - 1
    $ws.Range("C$r`:H$r").Value2 = $ws.Range("C$src`:H$src").Value2
}

# Write the brand-new first data row's sensor values
$ws.Range("C2").Value2 = -0.3135113716125494
$ws.Range("D2").Value2 = 1.868308603763581
$ws.Range("E2").Value2 = 2.184598565101624
$ws.Range("F2").Value2 = -0.0467311926186084
$ws.Range("G2").Value2 = 0.0064140851609408
$ws.Range("H2").Value2 = -0.0200058370828628
